# Updated stats for Mar 27
# - I21 (3/26) and I22 (3/27) become hardcoded "actual" values instead of
#   forecast formulas, which ripples through the forecast chain (I23:I40).
# - A new forecast row 41 (4/15) is appended, continuing the B/E/A/H/I/J/K/L/M/N
#   formula patterns from row 40.
# - A new annotation "EASTER SUNDAY..." is added at P38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Turn I21 / I22 into hard-coded "actual" numbers (style copied from the
#    existing "actual" cell I20, then the forecast formula is overwritten).
# ---------------------------------------------------------------------------
$ws.Range("I20").Copy()
$ws.Range("I21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I21").Value = 85435

$ws.Range("I22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I22").Value = 104126
$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) New annotation next to the Easter Sunday forecast row.
# ---------------------------------------------------------------------------
$ws.Range("P38").Value2 = "EASTER SUNDAY - THE DAY TRUMP WANTS CORONAVIRUS TO BE DONE"

# ---------------------------------------------------------------------------
# 3) Extend the table one more day: row 41 (4/15 forecast).
#    Copy formats down from row 40 first, then write the formulas (with
#    references shifted by one row) explicitly so they evaluate immediately.
# ---------------------------------------------------------------------------
$ws.Range("A40:N40").Copy()
$ws.Range("A41:N41").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = 0

$ws.Range("A41").Formula = "=+A40+1"
$ws.Range("B41").Formula = "=ROUND(B40+(B40*`$E`$10),0)"
$ws.Range("C41").Formula = "=+B41/`$P`$15*1000000"
$ws.Range("E41").Formula = "=+E40"
$ws.Range("F41").Value2 = $ws.Range("F40").Value2
$ws.Range("H41").Formula = "=A41"
$ws.Range("I41").Formula = "=I40*(1+AVERAGE(M38:M40))"
$ws.Range("J41").Formula = "=IF(I41<=0,0, I41-B41)"
$ws.Range("K41").Formula = "=IF(I41<=0,0, I41/`$P`$15*1000000)"
$ws.Range("L41").Formula = "=IF(I41<=0,0, I41-I40)"
$ws.Range("M41").Formula = "=IF(I41<=0,0, L41/I40)"
$ws.Range("N41").Formula = "=IF(I41<=0,0,M41-E41)"

# ---------------------------------------------------------------------------
# 4) Cosmetic view tweaks: column N a hair wider, scroll/selection moved to
#    the newly-updated cell.
# ---------------------------------------------------------------------------
$ws.Columns.Item(14).ColumnWidth = 3.83

$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("I23").Select()
